$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Content edits -------------------------------------------------
# Fix the "hw4?" label to "hw4" now that it has a real grade.
$ws.Range("A5").Value() = "hw4"

# Assignment 3 (project part1) grade corrected from 84 to 100.
$ws.Range("H3").Value() = 100

# Homework 3 grade entered (96) - this cell previously had a
# "needs grading" yellow highlight; clear it now that it's graded.
$ws.Range("C4").Value() = 96
$ws.Range("C4").Interior.ColorIndex = -4142

# Homework 4 (hw4) grade entered (94).
$ws.Range("C5").Value() = 94

# New "Estimated" column header, highlighted like the other headers.
$ws.Range("K1").Value() = "Estimated"
$ws.Range("K1").Interior.Color() = 65535
$ws.Columns.Item(11).ColumnWidth = 9.6

# Highlight the totals that now reflect the estimated final grade.
$ws.Range("C7").Interior.Color() = 65535
$ws.Range("H7").Interior.Color() = 65535
$ws.Range("B11").Interior.Color() = 65535

# Move the active selection to H4, matching the saved view state.
[void]$ws.Range("H4").Select()
